# ---------------------------------------------------------------------------
# Applies the "fix some problems with other devices" edit to tablet.xlsx:
#  * Apple sheet: renames/re-prices a couple of existing rows and adds the
#    rest of the current iPad line-up.
#  * Samsung sheet: adds the rest of the current Galaxy Tab / View line-up.
#  * Re-fonts everything from Calibri/JetBrains Mono to Times New Roman 12,
#    enlarges the header rows and switches the active tab to "Samsung".
# ---------------------------------------------------------------------------

$xlNone        = -4142
$xlThin        = 2
$xlMedium      = -4138
$xlLeft        = -4131
$xlCenter      = -4108
$xlEdgeLeft    = 7
$xlEdgeTop     = 8
$xlEdgeBottom  = 9
$xlEdgeRight   = 10

$MONEY_FMT   = '"$"#,##0.00'
$GENERAL_FMT = 'General'

function Set-OneEdge {
    param($rng, $idx, $w)
    if ($w -eq $xlNone) {
        $rng.Borders.Item($idx).LineStyle = $xlNone
    } else {
        $rng.Borders.Item($idx).LineStyle = 1
        $rng.Borders.Item($idx).Weight = $w
    }
}

function Set-Edges {
    param($rng, $left, $top, $right, $bottom)
    Set-OneEdge -rng $rng -idx $xlEdgeLeft   -w $left
    Set-OneEdge -rng $rng -idx $xlEdgeTop    -w $top
    Set-OneEdge -rng $rng -idx $xlEdgeRight  -w $right
    Set-OneEdge -rng $rng -idx $xlEdgeBottom -w $bottom
}

function Set-CellFormat {
    param(
        $rng,
        [bool]$bold,
        [bool]$colored,
        [string]$numFmt,
        [string]$hAlign = $null,
        [string]$vAlign = $null,
        $left, $top, $right, $bottom
    )
    $rng.Font.Name = "Times New Roman"
    $rng.Font.Size = 12
    $rng.Font.Bold = $bold
    if ($colored) { $rng.Font.ThemeColor = 1 }
    $rng.NumberFormat = $numFmt
    if ($hAlign -eq "left") { $rng.HorizontalAlignment = $xlLeft }
    if ($vAlign -eq "center") { $rng.VerticalAlignment = $xlCenter }
    Set-Edges -rng $rng -left $left -top $top -right $right -bottom $bottom
}

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Apple
$ws2 = $wb.Worksheets.Item(2)   # Samsung

# ===========================================================================
# Apple sheet
# ===========================================================================

$ws1.Rows.Item(1).RowHeight = 78

# Header row
$a1 = $ws1.Range("A1")
$a1.Value = "Модель"
Set-CellFormat -rng $a1 -bold $true -colored $true -numFmt $MONEY_FMT -left $xlMedium -top $xlMedium -right $xlMedium -bottom $xlMedium

$b1 = $ws1.Range("B1")
$b1.Value = "Стоимость"
Set-CellFormat -rng $b1 -bold $true -colored $true -numFmt $GENERAL_FMT -hAlign "left" -left $xlMedium -top $xlMedium -right $xlMedium -bottom $xlMedium

$ws1.Rows.Item(2).RowHeight = 15.75
$ws1.Rows.Item(3).RowHeight = 15.75
$ws1.Rows.Item(4).RowHeight = 15.75

# Row 2 : iPad Pro 2022 / 100000
$a2 = $ws1.Range("A2")
$a2.Value = "iPad Pro 2022"
Set-CellFormat -rng $a2 -bold $false -colored $false -numFmt $GENERAL_FMT -left $xlNone -top $xlMedium -right $xlThin -bottom $xlThin

$b2 = $ws1.Range("B2")
$b2.Value = 100000
Set-CellFormat -rng $b2 -bold $false -colored $true -numFmt $GENERAL_FMT -hAlign "left" -vAlign "center" -left $xlThin -top $xlNone -right $xlThin -bottom $xlThin

# Row 3 : iPad 2021 / 40000
$a3 = $ws1.Range("A3")
$a3.Value = "iPad 2021"
Set-CellFormat -rng $a3 -bold $false -colored $false -numFmt $GENERAL_FMT -left $xlNone -top $xlThin -right $xlThin -bottom $xlThin

$b3 = $ws1.Range("B3")
$b3.Value = 40000
Set-CellFormat -rng $b3 -bold $false -colored $true -numFmt $GENERAL_FMT -hAlign "left" -left $xlThin -top $xlThin -right $xlThin -bottom $xlThin

# Row 4 : iPad Air 2022 / 40000  (new row)
$a4 = $ws1.Range("A4")
$a4.Value = "iPad Air 2022"
Set-CellFormat -rng $a4 -bold $false -colored $false -numFmt $GENERAL_FMT -left $xlNone -top $xlNone -right $xlThin -bottom $xlThin

$b4 = $ws1.Range("B4")
$b4.Value = 40000
Set-CellFormat -rng $b4 -bold $false -colored $true -numFmt $GENERAL_FMT -hAlign "left" -left $xlThin -top $xlThin -right $xlThin -bottom $xlThin

# Rows 5-10 : remaining iPad models / 40000 (new rows)
$appleRows = @(
    @{ row = 5;  name = "iPad mini 2021" },
    @{ row = 6;  name = "iPad 2021" },
    @{ row = 7;  name = "iPad Pro 12,9 2021" },
    @{ row = 8;  name = "iPad Pro 11 2021" },
    @{ row = 9;  name = "iPad Air 2020" },
    @{ row = 10; name = "iPad 2020" }
)

foreach ($item in $appleRows) {
    $r = $item.row
    $ws1.Rows.Item($r).RowHeight = 15.75

    $aCell = $ws1.Cells.Item($r, 1)
    $aCell.Value = $item.name
    Set-CellFormat -rng $aCell -bold $false -colored $true -numFmt $MONEY_FMT -left $xlThin -top $xlThin -right $xlThin -bottom $xlThin

    $bCell = $ws1.Cells.Item($r, 2)
    $bCell.Value = 40000
    Set-CellFormat -rng $bCell -bold $false -colored $true -numFmt $GENERAL_FMT -hAlign "left" -left $xlThin -top $xlThin -right $xlThin -bottom $xlThin
}

# ===========================================================================
# Samsung sheet
# ===========================================================================

$ws2.Rows.Item(1).RowHeight = 62.25

$a1b = $ws2.Range("A1")
$a1b.Value = "Модель"
Set-CellFormat -rng $a1b -bold $true -colored $true -numFmt $GENERAL_FMT -hAlign "left" -left $xlThin -top $xlThin -right $xlThin -bottom $xlThin

$b1b = $ws2.Range("B1")
$b1b.Value = "Стоимость"
Set-CellFormat -rng $b1b -bold $true -colored $true -numFmt $GENERAL_FMT -hAlign "left" -left $xlThin -top $xlThin -right $xlThin -bottom $xlThin

$samsungRows = @(
    @{ row = 2;  name = "Samsung Galaxy Tab S7";         price = 15000; styleA = "header" },
    @{ row = 3;  name = "Samsung Galaxy Tab S6 10.5";    price = 17000; styleA = "header" },
    @{ row = 4;  name = "Samsung Galaxy Tab S5e 10.5";   price = 17000; styleA = "plain" },
    @{ row = 5;  name = "Samsung Galaxy Tab S4 10.5";    price = 17000; styleA = "plain" },
    @{ row = 6;  name = "Samsung Galaxy Tab S3 9.7";     price = 17000; styleA = "plain" },
    @{ row = 7;  name = "Samsung Galaxy Tab S2 9.7";     price = 17000; styleA = "plain" },
    @{ row = 8;  name = "Samsung Galaxy Tab S2 8.0";     price = 17000; styleA = "plain" },
    @{ row = 9;  name = "Samsung Galaxy Tab A 8.0";      price = 17000; styleA = "plain" },
    @{ row = 10; name = "Samsung Galaxy Tab A 10.1";     price = 17000; styleA = "plain" },
    @{ row = 11; name = "Samsung Galaxy Tab A 8.0";      price = 17000; styleA = "plain" },
    @{ row = 12; name = "Samsung Galaxy Tab A 10.5";     price = 17000; styleA = "plain" },
    @{ row = 13; name = "Samsung Galaxy Tab E 9.6";      price = 17000; styleA = "plain" },
    @{ row = 14; name = "Samsung Galaxy View 2";         price = 17000; styleA = "plain" }
)

foreach ($item in $samsungRows) {
    $r = $item.row
    $ws2.Rows.Item($r).RowHeight = 15.75

    $aCell = $ws2.Cells.Item($r, 1)
    $aCell.Value = $item.name
    if ($item.styleA -eq "header") {
        Set-CellFormat -rng $aCell -bold $false -colored $false -numFmt $GENERAL_FMT -left $xlNone -top $xlThin -right $xlThin -bottom $xlThin
    } else {
        Set-CellFormat -rng $aCell -bold $false -colored $true -numFmt $GENERAL_FMT -hAlign "left" -left $xlThin -top $xlThin -right $xlThin -bottom $xlThin
    }

    $bCell = $ws2.Cells.Item($r, 2)
    $bCell.Value = $item.price
    Set-CellFormat -rng $bCell -bold $false -colored $true -numFmt $GENERAL_FMT -hAlign "left" -left $xlThin -top $xlThin -right $xlThin -bottom $xlThin
}

# ===========================================================================
# Views / selection / active tab
# ===========================================================================

$ws1.Activate()
$ws1.Range("E7").Select()

$ws2.Activate()
$ws2.Range("C2").Select()

Write-Output "done"
